# Apply "Add data for 2022-06-14" update:
# - Rename sheet from "Through 2022-06-05" to "Through 2022-06-06"
# - Update header cell I1 text from "2022 (through 06-05)" to "2022 (through 06-06)"
# - Update June row (row 7) value in column I: 14 -> 18
# - Update Total row (row 14) value in column I: 678 -> 682

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-06-06"

$ws.Range("I1").Value = "2022 (through 06-06)"

$ws.Range("I7").Value = 18
$ws.Range("I14").Value = 682
